$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.931.92'
$ws.Range("E2").Value = '  +4.47%  '
$ws.Range("D3").Value = '1.880.25'
$ws.Range("E3").Value = '  +3.48%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '279.04'
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5303'
$ws.Range("E7").Value = '  +3.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3464'
$ws.Range("E8").Value = '  -1.52%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06969'
$ws.Range("E9").Value = '  +4.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.20'
$ws.Range("E10").Value = '  +1.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8094'
$ws.Range("E11").Value = '  -2.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07829'
$ws.Range("E12").Value = '  -1.00%  '
$ws.Range("D13").Value = '1.870.60'
$ws.Range("E13").Value = '  +2.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.202'
$ws.Range("E14").Value = '  +2.40%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.65'
$ws.Range("E15").Value = '  +3.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.63'
$ws.Range("E16").Value = '  +3.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9999'
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008073'
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").Value = '26.972.99'
$ws.Range("E20").Value = '  +4.40%  '
$ws.Range("D21").Value = '2.106.66'
$ws.Range("E21").Value = '  +2.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.763'
$ws.Range("E22").Value = '  +0.75%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.07'
$ws.Range("E23").Value = '  +0.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.201'
$ws.Range("E24").Value = '  +1.89%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.373'
$ws.Range("E25").Value = '  +8.63%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '147.00'
$ws.Range("E26").Value = '  +3.37%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.38'
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.664'
$ws.Range("E28").Value = '  -0.12%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.71'
$ws.Range("E29").Value = '  +3.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.381'
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.337'
$ws.Range("E31").Value = '  +2.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08907'
$ws.Range("E32").Value = '  +0.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04955'
$ws.Range("E33").Value = '  +1.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.184'
$ws.Range("E34").Value = '  +4.70%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7366'
$ws.Range("E35").Value = '  +0.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.887'
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.300'
$ws.Range("E37").Value = '  +4.65%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.397'
$ws.Range("E38").Value = '  +2.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01856'
$ws.Range("E39").Value = '  +0.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5172'
$ws.Range("E40").Value = '  -0.88%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9626'
$ws.Range("E41").Value = '  +0.92%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '116.48'
$ws.Range("E42").Value = '  +4.50%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.204'
$ws.Range("E43").Value = '  -0.02%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.130'
$ws.Range("E44").Value = '  +1.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9995'
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4522'
$ws.Range("E46").Value = '  -1.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1350'
$ws.Range("E47").Value = '  -1.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.369'
$ws.Range("E48").Value = '  +1.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.33'
$ws.Range("E49").Value = '  -0.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05954'
$ws.Range("E50").Value = '  +2.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.505'
$ws.Range("E51").Value = '  +0.17%  '
